# Update scripts with new TPM values.
# - Drop the old "Neutrophils" sender rows (rows 5-7 in the original sheet).
# - Re-point the remaining rows at the Nppc/Npr3 ligand-receptor pair
#   (columns B/C), keep/adjust the Target cluster in column D, and refresh
#   the numeric NATMI metrics (columns E:T) with the new TPM-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three trailing rows (old "Neutrophils" sender block) - this
# shrinks the used range from A1:T7 down to A1:T4.
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2: MuSCs -> ECs (via Nppc/Npr3)
$ws.Range("B2").Value = "Nppc"
$ws.Range("C2").Value = "Npr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.1248335
$ws.Range("H2").Value = 0.249667
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.2561565
$ws.Range("N2").Value = 0.512313
$ws.Range("O2").Value = 0.1741229440611523
$ws.Range("P2").Value = 0.1287532712126501
$ws.Range("Q2").Value = 0.03197691244275
$ws.Range("R2").Value = 0.127907649771
$ws.Range("S2").Value = 0.1741229440611523
$ws.Range("T2").Value = 0.1287532712126501

# Row 3: MuSCs -> FAPs (via Nppc/Npr3)
$ws.Range("B3").Value = "Nppc"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.1248335
$ws.Range("H3").Value = 0.249667
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.036780333333333
$ws.Range("N3").Value = 3.110341
$ws.Range("O3").Value = 0.7047537110504822
$ws.Range("P3").Value = 0.781683420754159
$ws.Range("Q3").Value = 0.1294249177411667
$ws.Range("R3").Value = 0.776549506447
$ws.Range("S3").Value = 0.7047537110504822
$ws.Range("T3").Value = 0.781683420754159

# Row 4: MuSCs -> MuSCs (via Nppc/Npr3)
$ws.Range("B4").Value = "Nppc"
$ws.Range("C4").Value = "Npr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.1248335
$ws.Range("H4").Value = 0.249667
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.1781875
$ws.Range("N4").Value = 0.356375
$ws.Range("O4").Value = 0.1211233448883654
$ws.Range("P4").Value = 0.08956330803319101
$ws.Range("Q4").Value = 0.02224376928125
$ws.Range("R4").Value = 0.088975077125
$ws.Range("S4").Value = 0.1211233448883654
$ws.Range("T4").Value = 0.08956330803319101
